# Add a new "tc048" worksheet right after "AddTest", seeded from the
# AddTest template (Epic/Feature/Requirement) plus new Tcname/
# Tcdescription/priority/QA columns, then make it the active/selected tab.

$wb = $excel.ActiveWorkbook

$addTest = $wb.Worksheets.Item("AddTest")
$tc048 = $wb.Worksheets.Add($null, $addTest)
$tc048.Name = "tc048"

$headers = @("Epic", "Feature", "Requirement", "Tcname", "Tcdescription", "priority", "QA")
$values  = @("Epic Mohit", "Mohit Feature", "RQ-489", "Unit testing ", "work", "Low", "Mohit Aman")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $tc048.Cells.Item(1, $i + 1).Value = $headers[$i]
    $tc048.Cells.Item(2, $i + 1).Value = $values[$i]
}

$headerRange = $tc048.Range("A1:G2")
$headerRange.WrapText = $true
$tc048.Rows.Item(1).RowHeight = 29
$tc048.Rows.Item(2).RowHeight = 29

# Clear AddTest's old single-cell selection, use a range selection instead
$addTest.Range("A1:C2").Select() | Out-Null

# tc048 becomes the selected/active sheet, with G8 selected
$tc048.Activate() | Out-Null
$tc048.Range("G8").Select() | Out-Null
